$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cell, $val) {
    $origStyle = $cell.Style
    $cell.NumberFormat = "@"
    $cell.Value = $val
    $cell.Style = $origStyle
}

# Row 2
$ws.Range("D2").Value = '45.267.59'
$ws.Range("E2").Value = '  -2.96%  '

# Row 3
$ws.Range("D3").Value = '2.386.92'
$ws.Range("E3").Value = '  +5.01%  '

# Row 4
Set-TextValue $ws.Range("D4") '0.998'
$ws.Range("E4").Value = '  -0.25%  '

# Row 5
Set-TextValue $ws.Range("D5") '294.08'
$ws.Range("E5").Value = '  -2.66%  '

# Row 6
Set-TextValue $ws.Range("D6") '94.62'
$ws.Range("E6").Value = '  -5.73%  '

# Row 7
Set-TextValue $ws.Range("D7") '0.558'
$ws.Range("E7").Value = '  -0.73%  '

# Row 8
Set-TextValue $ws.Range("D8") '0.999'
$ws.Range("E8").Value = '  -0.07%  '

# Row 9
Set-TextValue $ws.Range("D9") '0.500'
$ws.Range("E9").Value = '  -3.30%  '

# Row 10
Set-TextValue $ws.Range("D10") '34.10'
$ws.Range("E10").Value = '  -5.25%  '

# Row 11
Set-TextValue $ws.Range("D11") '0.0776'
$ws.Range("E11").Value = '  -0.92%  '

# Row 12
Set-TextValue $ws.Range("D12") '6.95'
$ws.Range("E12").Value = '  -3.50%  '

# Row 13
Set-TextValue $ws.Range("D13") '0.104'
$ws.Range("E13").Value = '  +0.82%  '

# Row 14
$ws.Range("D14").Value = '2.742.96'
$ws.Range("E14").Value = '  +4.66%  '

# Row 15
$ws.Range("D15").Value = '2.377.90'
$ws.Range("E15").Value = '  +4.73%  '

# Row 16
Set-TextValue $ws.Range("D16") '13.98'
$ws.Range("E16").Value = '  +2.36%  '

# Row 17
Set-TextValue $ws.Range("D17") '0.824'
$ws.Range("E17").Value = '  +2.93%  '

# Row 18
$ws.Range("D18").Value = '45.180.59'
$ws.Range("E18").Value = '  -3.12%  '

# Row 19
Set-TextValue $ws.Range("D19") '12.43'
$ws.Range("E19").Value = '  -4.12%  '

# Row 20
$ws.Range("D20").Value = '0.0₃0929'
$ws.Range("E20").Value = '  +0.06%  '

# Row 21
Set-TextValue $ws.Range("D21") '6.08'
$ws.Range("E21").Value = '  +2.55%  '

# Row 22
Set-TextValue $ws.Range("D22") '66.23'
$ws.Range("E22").Value = '  +1.44%  '

# Row 23
Set-TextValue $ws.Range("D23") '238.67'
$ws.Range("E23").Value = '  -3.51%  '

# Row 24
Set-TextValue $ws.Range("D24") '2.75'
$ws.Range("E24").Value = '  -3.52%  '

# Row 25
Set-TextValue $ws.Range("D25") '1.00'
$ws.Range("E25").Value = '  +0.04%  '

# Row 26
Set-TextValue $ws.Range("D26") '1.88'
$ws.Range("E26").Value = '  -0.50%  '

# Row 27
$ws.Range("B27").Value = 'InjectiveProtocol'
$ws.Range("C27").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
Set-TextValue $ws.Range("D27") '38.22'
$ws.Range("E27").Value = '  -10.33%  '

# Row 28
$ws.Range("B28").Value = 'Toncoin'
$ws.Range("C28").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
Set-TextValue $ws.Range("D28") '2.20'
$ws.Range("E28").Value = '  -1.75%  '

# Row 29
Set-TextValue $ws.Range("D29") '9.54'
$ws.Range("E29").Value = '  -2.05%  '

# Row 30
Set-TextValue $ws.Range("D30") '3.81'
$ws.Range("E30").Value = '  +15.85%  '

# Row 31
Set-TextValue $ws.Range("D31") '21.00'
$ws.Range("E31").Value = '  +5.74%  '

# Row 32
Set-TextValue $ws.Range("D32") '2.72'
$ws.Range("E32").Value = '  -2.75%  '

# Row 33
Set-TextValue $ws.Range("D33") '146.93'
$ws.Range("E33").Value = '  +0.08%  '

# Row 34
Set-TextValue $ws.Range("D34") '5.39'
$ws.Range("E34").Value = '  -2.15%  '

# Row 35
Set-TextValue $ws.Range("D35") '0.0759'
$ws.Range("E35").Value = '  -2.42%  '

# Row 36
$ws.Range("B36").Value = 'ARBITRUM'
$ws.Range("C36").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
Set-TextValue $ws.Range("D36") '1.97'
$ws.Range("E36").Value = '  +13.37%  '

# Row 37
$ws.Range("B37").Value = 'Kaspa'
$ws.Range("C37").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
Set-TextValue $ws.Range("D37") '0.112'
$ws.Range("E37").Value = '  -2.85%  '

# Row 38
$ws.Range("E38").Value = '  -1.45%  '

# Row 39
Set-TextValue $ws.Range("D39") '14.79'
$ws.Range("E39").Value = '  -7.50%  '

# Row 40
Set-TextValue $ws.Range("D40") '3.75'
$ws.Range("E40").Value = '  -3.96%  '

# Row 41
Set-TextValue $ws.Range("D41") '0.0295'
$ws.Range("E41").Value = '  -0.85%  '

# Row 42
Set-TextValue $ws.Range("D42") '3.18'
$ws.Range("E42").Value = '  -2.77%  '

# Row 43
$ws.Range("D43").Value = '1.936.52'
$ws.Range("E43").Value = '  +6.44%  '

# Row 44
Set-TextValue $ws.Range("D44") '0.996'
$ws.Range("E44").Value = '  -0.34%  '

# Row 45
Set-TextValue $ws.Range("D45") '89.76'
$ws.Range("E45").Value = '  -0.09%  '

# Row 46
Set-TextValue $ws.Range("D46") '1.72'
$ws.Range("E46").Value = '  -12.84%  '

# Row 47
Set-TextValue $ws.Range("D47") '8.48'
$ws.Range("E47").Value = '  +8.49%  '

# Row 48
Set-TextValue $ws.Range("D48") '15.14'
$ws.Range("E48").Value = '  +18.00%  '

# Row 49
Set-TextValue $ws.Range("D49") '99.69'
$ws.Range("E49").Value = '  +5.32%  '

# Row 50
$ws.Range("D50").Value = '2.613.89'
$ws.Range("E50").Value = '  +4.61%  '

# Row 51
Set-TextValue $ws.Range("D51") '0.182'
$ws.Range("E51").Value = '  -4.34%  '
